$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data rows (row, date-serial, B, C, D, E) appended after the previous
# last row (238), extending the table through 13/09/2021 (row 258).
$data = @(
    @(239, 44433, 137, 42, 47, 11),
    @(240, 44434, 127, 38, 47, 10),
    @(241, 44435, 127, 35, 47, 11),
    @(242, 44436, 127, 32, 47, 12),
    @(243, 44437, 127, 32, 47, 11),
    @(244, 44438, 127, 30, 47, 10),
    @(245, 44439, 127, 27, 47, 9),
    @(246, 44440, 127, 25, 47, 7),
    @(247, 44441, 127, 24, 47, 6),
    @(248, 44442, 127, 24, 47, 9),
    @(249, 44443, 127, 25, 48, 9),
    @(250, 44444, 127, 24, 48, 9),
    @(251, 44445, 127, 24, 48, 9),
    @(252, 44446, 127, 24, 47, 8),
    @(253, 44447, 127, 26, 38, 12),
    @(254, 44448, 127, 25, 38, 15),
    @(255, 44449, 127, 23, 38, 11),
    @(256, 44450, 127, 24, 38, 10),
    @(257, 44451, 127, 23, 38, 9),
    @(258, 44452, 127, 23, 38, 8)
)

$lastRow = 238
$firstNewRow = 239
$lastNewRow = 258

# Copy the formatting (number formats/styles) of the last existing row down
# across all the new rows in one shot, so the new date cells reuse the same
# date style (and F/G keep the percentage style) instead of minting new ones.
$ws.Range("A$lastRow`:G$lastRow").Copy() | Out-Null
$ws.Range("A$firstNewRow`:G$lastNewRow").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
[void]($excel.CutCopyMode = $false)

foreach ($row in $data) {
    $r  = $row[0]
    $dt = $row[1]
    $b  = $row[2]
    $c  = $row[3]
    $d  = $row[4]
    $e  = $row[5]

    $ws.Cells.Item($r, 1).Value = $dt
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e

    $ws.Cells.Item($r, 6).Formula = "=C$r/B$r"
    $ws.Cells.Item($r, 7).Formula = "=E$r/D$r"
}

# Update the sheet view to match the new scroll/selection position.
$ws.Range("D243").Select() | Out-Null
